$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell F1 with the same style as the other header cells (copy format from E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
$ws.Cells.Item(1, 6).Value = "time_taken"

# Fill in the time_taken values for each data row
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:41:17.385139"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:41:17.385149"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:41:17.385152"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:41:17.385155"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:41:17.385158"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:41:17.385160"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:41:17.385163"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:41:17.385165"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:41:17.385168"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:41:17.385170"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:41:17.385173"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:41:17.385175"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:41:17.385178"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:41:17.385180"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:41:17.385183"
